# [Fonds de solidarite] Add 2020-12-09 data
#
# Updates the "nombre_aides" (col C) and "montant_total" (col D) figures for
# the rows whose aggregates shifted with the 2020-12-09 data refresh. The
# source sheet stores every value (including the numeric-looking ones) as
# plain text, so each write temporarily forces Text number formatting before
# assigning the value, then clears that formatting again so the cell keeps
# its original (default) style while its content remains text rather than
# being re-interpreted as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Value
    )

    $range = $ws.Range($CellRef)
    $range.NumberFormat = "@"
    $range.Value = $Value
    $range.ClearFormats()
}

# Auvergne-Rhône-Alpes / H - Transports et entreposage
Set-TextValue "C8" "771"
Set-TextValue "D8" "1866524.90"

# Auvergne-Rhône-Alpes / I - Hébergement et restauration
Set-TextValue "C9" "1364"
Set-TextValue "D9" "10197813.05"

# Centre-Val de Loire / R - Arts, spectacles et activités récréatives
Set-TextValue "C62" "75"
Set-TextValue "D62" "714460.00"

# Hauts-de-France / G - Commerce
Set-TextValue "C127" "337"
Set-TextValue "D127" "1115116.00"

# Hauts-de-France / I - Hébergement et restauration
Set-TextValue "C129" "1120"
Set-TextValue "D129" "8961082.69"

# Hauts-de-France / M - Activités spécialisées, scientifiques et techniques
Set-TextValue "C133" "162"
Set-TextValue "D133" "1442979.82"

# Île-de-France / H - Transports et entreposage
Set-TextValue "C145" "8228"
Set-TextValue "D145" "25800767.91"

# Île-de-France / I - Hébergement et restauration
Set-TextValue "C146" "5005"
Set-TextValue "D146" "32769305.85"

# Île-de-France / N - Activités de services administratifs et de soutien
Set-TextValue "C151" "849"
Set-TextValue "D151" "3233731.22"

# Pays de la Loire / I - Hébergement et restauration
Set-TextValue "C244" "475"
Set-TextValue "D244" "3464962.16"

# Provence-Alpes-Côte d'Azur / H - Transports et entreposage
Set-TextValue "C259" "604"
Set-TextValue "D259" "1513025.18"
